$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$ws = $wb.Worksheets.Item("DPbES")

# ---------------------------------------------------------------------------
# 1. DPbES sheet: rename existing source labels that were split into more
#    specific categories ("coal" -> "hard coal", "wind" -> "onshore wind"),
#    then append five new source rows (13-17), in the same order the
#    original author must have touched them so the shared-string table
#    comes out in the same append order as the authoritative diff.
# ---------------------------------------------------------------------------

# New row 13: lignite (copies row 2 "hard coal" values/formulas)
$ws.Range("A13").Value = "lignite"
$ws.Range("B13").Formula = "=B2"
$ws.Range("C13:AK13").Formula = "=C2"

# Row 2 label: coal -> hard coal
$ws.Range("A2").Value = "hard coal"

# Row 6 label: wind -> onshore wind
$ws.Range("A6").Value = "onshore wind"

# New row 14: offshore wind (copies row 6 "onshore wind" values/formulas)
$ws.Range("A14").Value = "offshore wind"
$ws.Range("B14").Formula = "=B6"
$ws.Range("C14:AK14").Formula = "=C6"

# New row 15: crude oil (copies row 11 "petroleum" values/formulas)
$ws.Range("A15").Value = "crude oil"
$ws.Range("B15").Formula = "=B11"
$ws.Range("C15:AK15").Formula = "=C11"

# New row 16: heavy or residual fuel oil (copies row 11 "petroleum" values/formulas)
$ws.Range("A16").Value = "heavy or residual fuel oil"
$ws.Range("B16").Formula = "=B11"
$ws.Range("C16:AK16").Formula = "=C11"

# New row 17: municipal solid waste (copies row 9 "biomass" values/formulas)
$ws.Range("A17").Value = "municipal solid waste"
$ws.Range("B17").Formula = "=B9"
$ws.Range("C17:AK17").Formula = "=C9"

# New header cell A1: "Dispatch Priority (dimensionless)", bold + wrap text
$ws.Range("A1").Value = "Dispatch Priority (dimensionless)"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 30

# Widen column A slightly to fit the new header text
$ws.Columns.Item(1).ColumnWidth = 23

# ---------------------------------------------------------------------------
# 2. Tab selection moves from DPbES back to About.
# ---------------------------------------------------------------------------
$wsAbout.Select()
